$d = $word.ActiveDocument

# Locate the sole occurrence of "intr_ack" in the running text (not the
# table's "intr_en" register-name abbreviation, which already reads
# "intr_en" before this edit and must stay untouched).
$rng = $d.Content
$find = $rng.Find
$found = $find.Execute("intr_ack", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $wordStart = $rng.Start
    $wordEnd = $rng.End

    # Replace the trailing "ack" with "en" -> "intr_ack" becomes "intr_en"
    $ackRange = $d.Range($wordStart + 5, $wordEnd)
    $ackRange.Text = "en"

    # Word drops an invisible "_GoBack" bookmark at the location of the
    # most recent edit; mirror that here (zero-length, right after the
    # freshly typed "en").
    $editPoint = $d.Range($wordStart + 7, $wordStart + 7)
    $d.Bookmarks.Add("_GoBack", $editPoint)
}
